# Applies the midterm notice edits described by the commit:
#   "solved issues with document generation"
#
#  - MSICE -> None                              (4 occurrences)
#  - 2081  -> 2029                               (2 occurrences)
#  - 12 Mangsir   -> 14 mangsir,2080
#  - 22nd Mangsir -> 14 mangsir,2080
#  - 1:45 PM -> 1                                (2 occurrences)
#  - 5 minutes -> 5 mins
#  - Prof Dr. Baibhav  Singh -> Mr. user  user

$d = $word.ActiveDocument

# Word "Find & Replace" constants
$wdReplaceOne = 1
$wdFindContinue = 1

function Replace-InRange($range, $find, $replace) {
    $range.Find.Execute($find, $true, $true, $false, $false, $false,
                         $true, $wdFindContinue, $false, $replace, $wdReplaceOne) | Out-Null
}

# --- Heading line: "... Department of Electronics and Computer Engineering" / program name
Replace-InRange $d.Paragraphs(4).Range "MSICE" "None"

# --- "NOTICE FOR MID-TERM THESIS DEFENSE (2081 BATCH)"
Replace-InRange $d.Paragraphs(7).Range "2081" "2029"

# --- "Date: - 12 Mangsir"
Replace-InRange $d.Paragraphs(9).Range "12 Mangsir" "14 mangsir,2080"

# --- Big body paragraph: program, batch year, date, time, program again
$body = $d.Paragraphs(12).Range
Replace-InRange $body "MSICE" "None"
Replace-InRange $body "2081" "2029"
Replace-InRange $body "22nd Mangsir" "14 mangsir,2080"
Replace-InRange $body "1:45 PM" "1"
Replace-InRange $body "MSICE" "None"

# --- "Presentation Time: - 5 minutes (maximum)"
Replace-InRange $d.Paragraphs(14).Range "5 minutes" "5 mins"

# --- "Time: - 1:45 PM (Sharp)"
Replace-InRange $d.Paragraphs(15).Range "1:45 PM" "1"

# --- Signature table: coordinator name and program
$tbl = $d.Tables(1)
Replace-InRange $tbl.Cell(2, 1).Range "Prof Dr. Baibhav  Singh" "Mr. user  user"
Replace-InRange $tbl.Cell(4, 1).Range "MSICE" "None"
